# Convert bsecode (column D) cells that were stored as inline strings back
# into real numbers on the "day", "week", "month" and "quarter" sheets, and
# append the five new "week" rows (57-61) that were split out of stock.yaml.

$wb = $excel.ActiveWorkbook

# --- "day" sheet: D22:D49 -> numeric -------------------------------------
$ws = $wb.Worksheets.Item("day")
$dayCodes = [ordered]@{
    22 = 541729; 23 = 500510; 24 = 512599; 25 = 500325; 26 = 517354;
    27 = 532187; 28 = 532921; 29 = 532215; 30 = 531344; 31 = 532174;
    32 = 500112; 33 = 500440; 34 = 500547; 35 = 500104; 36 = 532955;
    37 = 532810; 38 = 533278; 39 = 500400; 40 = 532555; 41 = 534816;
    42 = 532898; 43 = 500049; 44 = 532134; 45 = 500312; 46 = 540065;
    47 = 532155; 48 = 532461; 49 = 532483
}
foreach ($r in $dayCodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $dayCodes[$r]
}

# --- "week" sheet: D52:D56 -> numeric -------------------------------------
$ws = $wb.Worksheets.Item("week")
$weekCodes = [ordered]@{
    52 = 532830; 53 = 532296; 54 = 532400; 55 = 532482; 56 = 500049
}
foreach ($r in $weekCodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $weekCodes[$r]
}

# --- "week" sheet: append new rows 57-61 ----------------------------------
# NOTE: the bsecode column (D) keeps its leading apostrophe so Excel stores
# these numeric-looking codes as text (inline string), matching the newly
# appended rows broken out of stock.yaml (unlike the pre-existing D22:D56
# cells above, which are genuine numbers).
$newWeekRows = @(
    @{row=57; sr=1; nsecode="ASTRAL";   name="Astral Poly Technik Limited";       bsecode="'532830"; per_chg=-2.33; close=2204.6;  volume=193812;    dt="19/06/2024 11:34:24"},
    @{row=58; sr=2; nsecode="GLENMARK"; name="Glenmark Pharmaceuticals Limited";  bsecode="'532296"; per_chg=-0.58; close=1235.7;  volume=811531;    dt="19/06/2024 11:34:24"},
    @{row=59; sr=3; nsecode="BSOFT";    name="Birlasoft Ltd";                     bsecode="'532400"; per_chg=1.03;  close=690.85;  volume=4391235;   dt="19/06/2024 11:34:24"},
    @{row=60; sr=4; nsecode="GRANULES"; name="Granules India Limited";            bsecode="'532482"; per_chg=-1.26; close=468.05;  volume=1010091;   dt="19/06/2024 11:34:24"},
    @{row=61; sr=5; nsecode="BEL";      name="Bharat Electronics Limited";        bsecode="'500049"; per_chg=-2.81; close=309.3;   volume=58698689;  dt="19/06/2024 11:34:24"}
)
foreach ($rowData in $newWeekRows) {
    $r = $rowData.row
    $ws.Cells.Item($r, 1).Value = $rowData.sr
    $ws.Cells.Item($r, 2).Value = $rowData.nsecode
    $ws.Cells.Item($r, 3).Value = $rowData.name
    $ws.Cells.Item($r, 4).Value = $rowData.bsecode
    $ws.Cells.Item($r, 5).Value = $rowData.per_chg
    $ws.Cells.Item($r, 6).Value = $rowData.close
    $ws.Cells.Item($r, 7).Value = $rowData.volume
    $ws.Cells.Item($r, 8).Value = "week"
    $ws.Cells.Item($r, 9).Value = $rowData.dt
}

# --- "month" sheet: D19:D34 -> numeric ------------------------------------
$ws = $wb.Worksheets.Item("month")
$monthCodes = [ordered]@{
    19 = 542652; 20 = 540005; 21 = 523642; 22 = 500114; 23 = 500790;
    24 = 532830; 25 = 532187; 26 = 532755; 27 = 500670; 28 = 540611;
    29 = 540133; 30 = 500875; 31 = 540065; 32 = 530005; 33 = 541153;
    34 = 532822
}
foreach ($r in $monthCodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $monthCodes[$r]
}

# --- "quarter" sheet: D28:D53 -> numeric ----------------------------------
$ws = $wb.Worksheets.Item("quarter")
$quarterCodes = [ordered]@{
    28 = 540699; 29 = 532777; 30 = 532488; 31 = 541729; 32 = 512599;
    33 = 533150; 34 = 500325; 35 = 539524; 36 = 500410; 37 = 542650;
    38 = 533398; 39 = 540716; 40 = 500575; 41 = 532921; 42 = 524804;
    43 = 500302; 44 = 532733; 45 = 532400; 46 = 540133; 47 = 540777;
    48 = 533155; 49 = 500085; 50 = 534816; 51 = 532523; 52 = 517334;
    53 = 532822
}
foreach ($r in $quarterCodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $quarterCodes[$r]
}
